$p = $ppt.ActivePresentation

# --- Slide 5: "Round-trip delay" bullet + rectangle diagram reposition ---
$s5 = $p.Slides.Item(5)

# Update the bullet text to mention the (t4 - t1) delay calculation
$contentShape5 = $s5.Shapes.Item(3)
$para5 = $contentShape5.TextFrame.TextRange.Paragraphs(6)
$para5.Runs(1).Text = "Round-trip delay (t4 - t1) metrics are notified when consecutive M number of probe messages have delay values exceed the configured thresholds"

# Move the "Rectangle 2" diagram box up slightly (y offset 893624 -> 819150 EMU)
$rect5 = $s5.Shapes.Item(4)
$rect5.Top = 819150 / 12700

# --- Slide 6: "One-way delay" bullet + content placeholder resize ---
$s6 = $p.Slides.Item(6)

# Grow the content placeholder's height (cy 2138362 -> 2286000 EMU)
$contentShape6 = $s6.Shapes.Item(3)
$contentShape6.Height = 2286000 / 12700

# Update the bullet text to mention the (t2 - t1) delay calculation
$para6 = $contentShape6.TextFrame.TextRange.Paragraphs(6)
$para6.Runs(1).Text = "One-way delay (t2 – t1) metrics are notified when consecutive M number of probe messages have delay values exceed the configured thresholds"
